$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header/date formatting from the Weekly Quantity sheet so the new
# sheet reuses the same style definitions (bold/centered header, date format).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Populate header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Populate data rows ---
$data = @(
    @(45480.99999999999, 2, 1.999999997338895, 2.000000002561412),
    @(45592.99999999999, 2, 1.999999997277968, 2.00000000256242),
    @(45634.99999999999, 2, 1.999999997543789, 2.000000002647953),
    @(45641.99999999999, 2, 1.99999999746279,  2.000000002812607),
    @(45648.99999999999, 2, 1.99999999737104,  2.000000002516354),
    @(45655.99999999999, 2, 1.999999997283658, 2.000000002711568),
    @(45662.99999999999, 2, 1.999999997307405, 2.000000002812351),
    @(45669.99999999999, 2, 1.999999997187781, 2.000000002870865),
    @(45676.99999999999, 2, 1.999999997094176, 2.000000003075768),
    @(45683.99999999999, 2, 1.999999996964093, 2.000000003028802),
    @(45690.99999999999, 2, 1.999999996735448, 2.000000003151815)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
